$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2131979695431472
$ws.Range("C2").Value = 0.5126903553299492
$ws.Range("J2").Value = 0.06768189509306261
$ws.Range("O2").Value = 0.001692047377326565
$ws.Range("P2").Value = 0.1404399323181049
$ws.Range("S2").Value = 0.06429780033840947
$ws.Range("B3").Value = 0.01238390092879257
$ws.Range("C3").Value = 0.04024767801857585
$ws.Range("J3").Value = 0.1269349845201238
$ws.Range("P3").Value = 0.7461300309597523
$ws.Range("S3").Value = 0.07430340557275542
$ws.Range("J4").Value = 0.1578947368421053
$ws.Range("P4").Value = 0.6973684210526315
$ws.Range("S4").Value = 0.1447368421052632
$ws.Range("B6").Value = 0.06930693069306931
$ws.Range("D6").Value = 0.009900990099009901
$ws.Range("E6").Value = 0.004950495049504951
$ws.Range("F6").Value = 0.07425742574257425
$ws.Range("J6").Value = 0.349009900990099
$ws.Range("O6").Value = 0.04207920792079208
$ws.Range("Q6").Value = 0.1163366336633663
$ws.Range("R6").Value = 0.08415841584158416
$ws.Range("S6").Value = 0.25
$ws.Range("B7").Value = 0.08787878787878788
$ws.Range("D7").Value = 0.02727272727272727
$ws.Range("E7").Value = 0.00303030303030303
$ws.Range("F7").Value = 0.06666666666666667
$ws.Range("J7").Value = 0.2303030303030303
$ws.Range("O7").Value = 0.02727272727272727
$ws.Range("Q7").Value = 0.1606060606060606
$ws.Range("R7").Value = 0.08484848484848485
$ws.Range("S7").Value = 0.3121212121212121
$ws.Range("B8").Value = 0.1070588235294118
$ws.Range("D8").Value = 0.01529411764705882
$ws.Range("E8").Value = 0.001176470588235294
$ws.Range("F8").Value = 0.05764705882352941
$ws.Range("J8").Value = 0.1905882352941176
$ws.Range("O8").Value = 0.02
$ws.Range("Q8").Value = 0.1682352941176471
$ws.Range("R8").Value = 0.1035294117647059
$ws.Range("S8").Value = 0.3364705882352941
$ws.Range("B9").Value = 0.119047619047619
$ws.Range("D9").Value = 0.02040816326530612
$ws.Range("F9").Value = 0.04761904761904762
$ws.Range("J9").Value = 0.1836734693877551
$ws.Range("O9").Value = 0.01360544217687075
$ws.Range("Q9").Value = 0.1598639455782313
$ws.Range("R9").Value = 0.1258503401360544
$ws.Range("S9").Value = 0.3299319727891156
$ws.Range("B10").Value = 0.09827255278310941
$ws.Range("D10").Value = 0.01880998080614204
$ws.Range("E10").Value = 0.0007677543186180423
$ws.Range("F10").Value = 0.06641074856046066
$ws.Range("J10").Value = 0.2652591170825336
$ws.Range("O10").Value = 0.02495201535508637
$ws.Range("Q10").Value = 0.1915547024952015
$ws.Range("R10").Value = 0.08061420345489444
$ws.Range("S10").Value = 0.2533589251439539
$ws.Range("G11").Value = 0.1043256997455471
$ws.Range("J11").Value = 0.08396946564885496
$ws.Range("K11").Value = 0.1246819338422392
$ws.Range("L11").Value = 0.6717557251908397
$ws.Range("S11").Value = 0.01526717557251908
$ws.Range("G12").Value = 0.7862318840579711
$ws.Range("J12").Value = 0.1739130434782609
$ws.Range("L12").Value = 0.03985507246376811
$ws.Range("F13").Value = 0.01388888888888889
$ws.Range("G13").Value = 0.6944444444444444
$ws.Range("J13").Value = 0.2361111111111111
$ws.Range("S13").Value = 0.05555555555555555
$ws.Range("F15").Value = 0.02956989247311828
$ws.Range("H15").Value = 0.1666666666666667
$ws.Range("I15").Value = 0.05376344086021505
$ws.Range("J15").Value = 0.3897849462365591
$ws.Range("K15").Value = 0.06182795698924731
$ws.Range("M15").Value = 0.01344086021505376
$ws.Range("N15").Value = 0.002688172043010753
$ws.Range("O15").Value = 0.0456989247311828
$ws.Range("S15").Value = 0.2365591397849462
$ws.Range("F16").Value = 0.02645502645502645
$ws.Range("H16").Value = 0.2037037037037037
$ws.Range("I16").Value = 0.06349206349206349
$ws.Range("J16").Value = 0.4100529100529101
$ws.Range("K16").Value = 0.1111111111111111
$ws.Range("M16").Value = 0.02116402116402116
$ws.Range("N16").Value = 0.002645502645502645
$ws.Range("O16").Value = 0.06084656084656084
$ws.Range("S16").Value = 0.1005291005291005
$ws.Range("F17").Value = 0.02319587628865979
$ws.Range("H17").Value = 0.1842783505154639
$ws.Range("I17").Value = 0.06572164948453608
$ws.Range("J17").Value = 0.4768041237113402
$ws.Range("K17").Value = 0.09664948453608248
$ws.Range("M17").Value = 0.01804123711340206
$ws.Range("N17").Value = 0.001288659793814433
$ws.Range("O17").Value = 0.07860824742268041
$ws.Range("S17").Value = 0.05541237113402062
$ws.Range("F18").Value = 0.0202020202020202
$ws.Range("H18").Value = 0.1843434343434343
$ws.Range("I18").Value = 0.1035353535353535
$ws.Range("J18").Value = 0.4924242424242424
$ws.Range("K18").Value = 0.0707070707070707
$ws.Range("M18").Value = 0.0202020202020202
$ws.Range("O18").Value = 0.04040404040404041
$ws.Range("S18").Value = 0.06818181818181818
$ws.Range("F19").Value = 0.01930036188178528
$ws.Range("H19").Value = 0.2273823884197829
$ws.Range("I19").Value = 0.07418576598311219
$ws.Range("J19").Value = 0.4083232810615199
$ws.Range("K19").Value = 0.09710494571773221
$ws.Range("M19").Value = 0.01990349819059107
$ws.Range("N19").Value = 0.0006031363088057901
$ws.Range("O19").Value = 0.05729794933655006
$ws.Range("S19").Value = 0.09589867310012062
